$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Taxonsorteringsordning) values for rows 2-19: add 4 to each value
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value2 = $cell.Value2 + 4
}
